# Update scripts with new TPM values for the Cx3cl1-Cx3cr1 LR-pair sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (row -> A..T). Columns A-D are cluster/gene-symbol labels
# (text, stored as shared strings); columns E-T are numeric TPM-derived stats.
# Row 10 from the previous export no longer exists with the refreshed TPM
# numbers, so it gets removed below.

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Cx3cl1"
$ws.Cells.Item(2, 3).Value = "Cx3cr1"
$ws.Cells.Item(2, 4).Value = "MuSCs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 7.361448666666667
$ws.Cells.Item(2, 8).Value = 22.084346
$ws.Cells.Item(2, 9).Value = 0.3809728075517136
$ws.Cells.Item(2, 10).Value = 0.3809728075517136
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.004344
$ws.Cells.Item(2, 14).Value = 0.013032
$ws.Cells.Item(2, 15).Value = [double]"8.384535974127607E-05"
$ws.Cells.Item(2, 16).Value = [double]"8.384535974127607E-05"
$ws.Cells.Item(2, 17).Value = 0.031978133008
$ws.Cells.Item(2, 18).Value = 0.287803197072
$ws.Cells.Item(2, 19).Value = [double]"3.194280210081737E-05"
$ws.Cells.Item(2, 20).Value = [double]"3.194280210081737E-05"

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Cx3cl1"
$ws.Cells.Item(3, 3).Value = "Cx3cr1"
$ws.Cells.Item(3, 4).Value = "Resolving-Mac"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 7.361448666666667
$ws.Cells.Item(3, 8).Value = 22.084346
$ws.Cells.Item(3, 9).Value = 0.3809728075517136
$ws.Cells.Item(3, 10).Value = 0.3809728075517136
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 51.805321
$ws.Cells.Item(3, 14).Value = 155.415963
$ws.Cells.Item(3, 15).Value = 0.9999161546402586
$ws.Cells.Item(3, 16).Value = 0.9999161546402586
$ws.Cells.Item(3, 17).Value = 381.3622112016887
$ws.Cells.Item(3, 18).Value = 3432.259900815198
$ws.Cells.Item(3, 19).Value = 0.3809408647496128
$ws.Cells.Item(3, 20).Value = 0.3809408647496128

$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Cx3cl1"
$ws.Cells.Item(4, 3).Value = "Cx3cr1"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 10.317205
$ws.Cells.Item(4, 8).Value = 30.951615
$ws.Cells.Item(4, 9).Value = 0.5339403605073807
$ws.Cells.Item(4, 10).Value = 0.5339403605073807
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.004344
$ws.Cells.Item(4, 14).Value = 0.013032
$ws.Cells.Item(4, 15).Value = [double]"8.384535974127607E-05"
$ws.Cells.Item(4, 16).Value = [double]"8.384535974127607E-05"
$ws.Cells.Item(4, 17).Value = 0.04481793852
$ws.Cells.Item(4, 18).Value = 0.40336144668
$ws.Cells.Item(4, 19).Value = [double]"4.476842160712797E-05"
$ws.Cells.Item(4, 20).Value = [double]"4.476842160712797E-05"

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Cx3cl1"
$ws.Cells.Item(5, 3).Value = "Cx3cr1"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 10.317205
$ws.Cells.Item(5, 8).Value = 30.951615
$ws.Cells.Item(5, 9).Value = 0.5339403605073807
$ws.Cells.Item(5, 10).Value = 0.5339403605073807
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 51.805321
$ws.Cells.Item(5, 14).Value = 155.415963
$ws.Cells.Item(5, 15).Value = 0.9999161546402586
$ws.Cells.Item(5, 16).Value = 0.9999161546402586
$ws.Cells.Item(5, 17).Value = 534.486116847805
$ws.Cells.Item(5, 18).Value = 4810.375051630244
$ws.Cells.Item(5, 19).Value = 0.5338955920857735
$ws.Cells.Item(5, 20).Value = 0.5338955920857735

$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 2).Value = "Cx3cl1"
$ws.Cells.Item(6, 3).Value = "Cx3cr1"
$ws.Cells.Item(6, 4).Value = "MuSCs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.634232333333333
$ws.Cells.Item(6, 8).Value = 4.902697
$ws.Cells.Item(6, 9).Value = 0.08457548349701474
$ws.Cells.Item(6, 10).Value = 0.08457548349701474
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.004344
$ws.Cells.Item(6, 14).Value = 0.013032
$ws.Cells.Item(6, 15).Value = [double]"8.384535974127607E-05"
$ws.Cells.Item(6, 16).Value = [double]"8.384535974127607E-05"
$ws.Cells.Item(6, 17).Value = 0.007099105255999999
$ws.Cells.Item(6, 18).Value = 0.063891947304
$ws.Cells.Item(6, 19).Value = [double]"7.091261839099558E-06"
$ws.Cells.Item(6, 20).Value = [double]"7.091261839099558E-06"

$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Cx3cl1"
$ws.Cells.Item(7, 3).Value = "Cx3cr1"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.634232333333333
$ws.Cells.Item(7, 8).Value = 4.902697
$ws.Cells.Item(7, 9).Value = 0.08457548349701474
$ws.Cells.Item(7, 10).Value = 0.08457548349701474
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 51.805321
$ws.Cells.Item(7, 14).Value = 155.415963
$ws.Cells.Item(7, 15).Value = 0.9999161546402586
$ws.Cells.Item(7, 16).Value = 0.9999161546402586
$ws.Cells.Item(7, 17).Value = 84.66193061691233
$ws.Cells.Item(7, 18).Value = 761.957375552211
$ws.Cells.Item(7, 19).Value = 0.08456839223517564
$ws.Cells.Item(7, 20).Value = 0.08456839223517564

$ws.Cells.Item(8, 1).Value = "Resolving-Mac"
$ws.Cells.Item(8, 2).Value = "Cx3cl1"
$ws.Cells.Item(8, 3).Value = "Cx3cr1"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.009880666666666668
$ws.Cells.Item(8, 8).Value = 0.029642
$ws.Cells.Item(8, 9).Value = 0.0005113484438908852
$ws.Cells.Item(8, 10).Value = 0.0005113484438908852
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.004344
$ws.Cells.Item(8, 14).Value = 0.013032
$ws.Cells.Item(8, 15).Value = [double]"8.384535974127607E-05"
$ws.Cells.Item(8, 16).Value = [double]"8.384535974127607E-05"
$ws.Cells.Item(8, 17).Value = [double]"4.2921616E-05"
$ws.Cells.Item(8, 18).Value = 0.000386294544
$ws.Cells.Item(8, 19).Value = [double]"4.287419423117299E-08"
$ws.Cells.Item(8, 20).Value = [double]"4.287419423117299E-08"

$ws.Cells.Item(9, 1).Value = "Resolving-Mac"
$ws.Cells.Item(9, 2).Value = "Cx3cl1"
$ws.Cells.Item(9, 3).Value = "Cx3cr1"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.009880666666666668
$ws.Cells.Item(9, 8).Value = 0.029642
$ws.Cells.Item(9, 9).Value = 0.0005113484438908852
$ws.Cells.Item(9, 10).Value = 0.0005113484438908852
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 51.805321
$ws.Cells.Item(9, 14).Value = 155.415963
$ws.Cells.Item(9, 15).Value = 0.9999161546402586
$ws.Cells.Item(9, 16).Value = 0.9999161546402586
$ws.Cells.Item(9, 17).Value = 0.5118711083606667
$ws.Cells.Item(9, 18).Value = 4.606839975246
$ws.Cells.Item(9, 19).Value = 0.000511305569696654
$ws.Cells.Item(9, 20).Value = 0.000511305569696654

# The old row 10 (MuSCs -> Resolving-Mac under the stale TPM numbers) is no
# longer part of the refreshed table, so remove it entirely and let the
# sheet's used range shrink to A1:T9.
$ws.Rows.Item(10).Delete()
